# Flash cards workbook update:
#  - delete the empty placeholder sheet "Varieties of designations" (sheetId 6)
#  - rename "European design. & varieties" (sheetId 4) to "Varieties of designations"
#  - clear the ad-hoc cell selections that were left over on several sheets,
#    restoring them to their natural/default state
#  - keep "Varieties of designations" (renamed sheet, ex "European design. & varieties")
#    as the active/selected tab

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the stray empty "Varieties of designations" sheet (sheetId 6 / rId10)
$wb.Worksheets.Item("Varieties of designations").Delete()

# 2. Rename "European design. & varieties" (sheetId 4 / rId5) now that the name is free
$euroSheet = $wb.Worksheets.Item("European design. & varieties")
$euroSheet.Name = "Varieties of designations"

# 3. Tidy up leftover selections on the sheets that no longer carry a meaningful one
$wb.Worksheets.Item("Wine tasting").Range("A1").Select()
$wb.Worksheets.Item("Nobles varieties").Range("A1").Select()
$wb.Worksheets.Item("Wine to discover (todo)").Range("A1").Select()

# 4. "Wine & naming convention": selection moves up one row, scroll resets to top
$wb.Worksheets.Item("Wine & naming convention").Range("A19").Select()

# 5. Re-activate the renamed sheet so it stays the active/visible tab, selection reset to A1
$euroSheet.Activate()
$euroSheet.Range("A1").Select()
